$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hungarian Excel shows the built-in "Normal" cell style as "Normál".
$wb.Styles.Item(1).Name = "Normál"

# Set cells in the same order new shared-strings were introduced, so the
# resulting shared string table ordering matches the author's edit flow.
$ws.Range("C3").Value = "Felhasználókezelés - utánaolvasni, hogy hogyan működik - Angular és ASP.NET authentication"
$ws.Range("C5").Value = "Felhasználókezelés - felület létrehozás - bejelentkezés, profil, rendeléseim képernyő"
$ws.Range("C6").Value = "Shop filters - utánaolvasni, hogyan szokás elkészíteni, backend queryk? Megvalósítása"
$ws.Range("C7").Value = "Shop filters - webes kliensben megvalósítás - shopban filter felület, filter service?"
$ws.Range("C10").Value = "Fizetés - utánanézni, milyen lehetőségek vannak, hogyan lehet beépíteni őket, ki lehet-e próbálni őket"
$ws.Range("C11").Value = "Webes fizetés megvalósítása"
$ws.Range("C12").Value = "Kereső optimalizálás - utánanézni, milyen módszerek vannak rá, hogyan érdemes csinálni, implementálni"
$ws.Range("A16").Value = "nyár"
$ws.Range("C2").Value = "Feladatkiírás, projekt rendbeszedése, ütemterv részletesebb elkészítése"
$ws.Range("C4").Value = "Üzleti folyamat feltérképezés, Felhasználókezelés - backend megvalósítás, frontenden servicek megvalósítása, "
$ws.Range("A20").Value = "TODO: Android wireframe"

# Cells that reuse existing shared strings.
$ws.Range("C8").Value = "Admin felület (áruk hozzáadása, törlése, szerkesztése, rendelések nézése)"
$ws.Range("C9").Value = "Admin felület (áruk hozzáadása, törlése, szerkesztése, rendelések nézése)"
$ws.Range("C13").Value = "Android kliens"
$ws.Range("C14").Value = "Android kliens"
$ws.Range("C15").Value = "Android kliens"
$ws.Range("C16").Value = "Android kliens"
$ws.Range("A17").Value = "2. félév:"
$ws.Range("C17").Value = "Ajánló motor, tesztelés, deployment"

$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Columns.Item(1).ColumnWidth = 23.5
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(3).ColumnWidth = 108.83333333333333

$ws.Range("C19").Select()
